# Corrección todos los ingresos y TC
# Update the "missing_values" sheet with corrected income counts and
# their recalculated percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missing_values")

# --- Block 1: rows 7-9 (Monetary labor income) ---
$ws.Range("B7").Value = 36
$ws.Range("C7").Value = 0.29756984625557942

$ws.Range("B8").Value = 2197
$ws.Range("C8").Value = 18.160026450653

$ws.Range("B9").Value = 9700
$ws.Range("C9").Value = 80.178541907753342

# --- Block 2: rows 18-20 (Pensions and retirement benefits) ---
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 0.10193679918450561

$ws.Range("B19").Value = 171
$ws.Range("C19").Value = 2.90519877675841

$ws.Range("B20").Value = 5658
$ws.Range("C20").Value = 96.126401630988795

# --- Block 3: rows 38-40 (Monetary non-labor income) ---
$ws.Range("B38").Value = 86
$ws.Range("C38").Value = 1.4771556166265889

$ws.Range("B39").Value = 241
$ws.Range("C39").Value = 4.13947097217451

$ws.Range("B40").Value = 5577
$ws.Range("C40").Value = 95.791824115424248

$wb.Save()
